$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1727.75
$ws.Range("I28").Value = 526.26666
$ws.Range("K28").Value = 526.26666
$ws.Range("M28").Value = -41.26666

$ws.Range("H64").Value = 30306920
$ws.Range("J64").Value = 4474.5
$ws.Range("L64").Value = 4474.5
$ws.Range("N64").Value = -4970.5

$ws.Range("H67").Value = 30306920
$ws.Range("J67").Value = 4474.5
$ws.Range("L67").Value = 4474.5
$ws.Range("N67").Value = -6190.5

$ws.Range("H88").Value = 2457.2727
$ws.Range("I88").Value = 1783.25
$ws.Range("J88").Value = 2842.4285
$ws.Range("K88").Value = 1783.25
$ws.Range("L88").Value = 2842.4285
$ws.Range("M88").Value = -1377.25
$ws.Range("N88").Value = -3654.4285

$ws.Range("H91").Value = 2457.2727
$ws.Range("I91").Value = 1783.25
$ws.Range("J91").Value = 2842.4285
$ws.Range("K91").Value = 1783.25
$ws.Range("L91").Value = 2842.4285
$ws.Range("M91").Value = -379.25
$ws.Range("N91").Value = -5650.4285

$ws.Range("H98").Value = 740.4074000000001
$ws.Range("I98").Value = 760.0417
$ws.Range("K98").Value = 760.0417
$ws.Range("M98").Value = 737.9583

$ws.Range("H122").Value = 740.4074000000001
$ws.Range("I122").Value = 760.0417
$ws.Range("K122").Value = 2280.1251
$ws.Range("M122").Value = 169.8748999999998

$ws.Range("H129").Value = 22231722
$ws.Range("I129").Value = 3231
$ws.Range("J129").Value = 37050716
$ws.Range("K129").Value = 9693
$ws.Range("L129").Value = 111152148
$ws.Range("M129").Value = -4693
$ws.Range("N129").Value = -111162148

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5219.593
$ws.Range("I2").Value = 5101.5264
$ws.Range("K2").Value = 5101.5264
$ws.Range("M2").Value = -4988.5264

$ws.Range("H45").Value = 2065.5293
$ws.Range("I45").Value = 2186.75
$ws.Range("J45").Value = 1774.6
$ws.Range("K45").Value = 2186.75
$ws.Range("L45").Value = 1774.6
$ws.Range("M45").Value = -1809.75
$ws.Range("N45").Value = -2528.6

$ws.Range("H56").Value = 103332.664
$ws.Range("J56").Value = 103332.664
$ws.Range("L56").Value = 103332.664
$ws.Range("N56").Value = -104816.664

$ws.Range("H61").Value = 3398.8
$ws.Range("I61").Value = 1814.6
$ws.Range("J61").Value = 3926.8667
$ws.Range("K61").Value = 1814.6
$ws.Range("L61").Value = 3926.8667
$ws.Range("M61").Value = -1602.6
$ws.Range("N61").Value = -4350.8667

$ws.Range("H63").Value = 3257.7
$ws.Range("I63").Value = 1554
$ws.Range("K63").Value = 1554
$ws.Range("M63").Value = -868

$ws.Range("H66").Value = 3257.7
$ws.Range("I66").Value = 1554
$ws.Range("K66").Value = 7770
$ws.Range("M66").Value = -4338

$ws.Range("H74").Value = 2597.8462
$ws.Range("I74").Value = 1354.6
$ws.Range("J74").Value = 3374.875
$ws.Range("K74").Value = 1354.6
$ws.Range("L74").Value = 3374.875
$ws.Range("M74").Value = -480.5999999999999
$ws.Range("N74").Value = -5122.875

$ws.Range("H77").Value = 2597.8462
$ws.Range("I77").Value = 1354.6
$ws.Range("J77").Value = 3374.875
$ws.Range("K77").Value = 6773
$ws.Range("L77").Value = 16874.375
$ws.Range("M77").Value = -2405
$ws.Range("N77").Value = -25610.375

$ws.Range("H97").Value = 1804.2
$ws.Range("I97").Value = 1803.3334
$ws.Range("J97").Value = 1805.5
$ws.Range("K97").Value = 1803.3334
$ws.Range("L97").Value = 1805.5
$ws.Range("M97").Value = -1307.3334
$ws.Range("N97").Value = -2797.5

$ws.Range("H104").Value = 65000
$ws.Range("J104").Value = 65000
$ws.Range("L104").Value = 65000
$ws.Range("N104").Value = -71988

$ws.Range("H116").Value = 5219.593
$ws.Range("I116").Value = 5101.5264
$ws.Range("K116").Value = 5101.5264
$ws.Range("M116").Value = -2807.5264

$ws.Range("H122").Value = 2210.963
$ws.Range("I122").Value = 1330.7778
$ws.Range("J122").Value = 3971.3333
$ws.Range("K122").Value = 3992.3334
$ws.Range("L122").Value = 11913.9999
$ws.Range("M122").Value = -1542.3334
$ws.Range("N122").Value = -16813.9999

$ws.Range("H132").Value = 1560845.5
$ws.Range("I132").Value = 1894760.8
$ws.Range("J132").Value = 225184.75
$ws.Range("K132").Value = 5684282.4
$ws.Range("L132").Value = 675554.25
$ws.Range("M132").Value = -5681752.4
$ws.Range("N132").Value = -680614.25

$ws.Range("H136").Value = 3398.8
$ws.Range("I136").Value = 1814.6
$ws.Range("J136").Value = 3926.8667
$ws.Range("K136").Value = 5443.799999999999
$ws.Range("L136").Value = 11780.6001
$ws.Range("M136").Value = -2893.799999999999
$ws.Range("N136").Value = -16880.6001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5219.593
$ws.Range("I3").Value = 5101.5264
$ws.Range("K3").Value = 5101.5264
$ws.Range("M3").Value = -4987.5264

$ws.Range("H105").Value = 3124.6667
$ws.Range("I105").Value = 3471.2354
$ws.Range("K105").Value = 3471.2354
$ws.Range("M105").Value = -1724.2354

$ws.Range("H107").Value = 11115227
$ws.Range("J107").Value = 3149.7144
$ws.Range("L107").Value = 3149.7144
$ws.Range("N107").Value = -6989.7144

$ws.Range("H110").Value = 53732.668
$ws.Range("J110").Value = 53732.668
$ws.Range("L110").Value = 53732.668
$ws.Range("N110").Value = -61912.668

$ws.Range("H134").Value = 2986515
$ws.Range("I134").Value = 3972229
$ws.Range("K134").Value = 11916687
$ws.Range("M134").Value = -11914152

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1666.3334
$ws.Range("I107").Value = 1499.5
$ws.Range("K107").Value = 1499.5
$ws.Range("M107").Value = 420.5

$ws.Range("H132").Value = 4892.8057
$ws.Range("I132").Value = 3864.6553
$ws.Range("K132").Value = 11593.9659
$ws.Range("M132").Value = -9063.965899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1925.1428
$ws.Range("I14").Value = 1925.1428
$ws.Range("K14").Value = 5775.428400000001
$ws.Range("M14").Value = -5602.428400000001

$ws.Range("H133").Value = 3399.5
$ws.Range("I133").Value = 3399.5
$ws.Range("K133").Value = 10198.5
$ws.Range("M133").Value = -5138.5

$ws.Range("H139").Value = 2552.889
$ws.Range("I139").Value = 2139.4285
$ws.Range("K139").Value = 6418.2855
$ws.Range("M139").Value = -1278.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6129.317
$ws.Range("I102").Value = 4854.32
$ws.Range("K102").Value = 4854.32
$ws.Range("M102").Value = -3232.32

$ws.Range("H113").Value = 7970.8237
$ws.Range("I113").Value = 2100.8333
$ws.Range("J113").Value = 11172.637
$ws.Range("K113").Value = 2100.8333
$ws.Range("L113").Value = 11172.637
$ws.Range("M113").Value = 69.16670000000022
$ws.Range("N113").Value = -15512.637

$ws.Range("H122").Value = 7309.2
$ws.Range("I122").Value = 5482.1113
$ws.Range("J122").Value = 10049.833
$ws.Range("K122").Value = 16446.3339
$ws.Range("L122").Value = 30149.499
$ws.Range("M122").Value = -13996.3339
$ws.Range("N122").Value = -35049.499

$ws.Range("H126").Value = 22733354
$ws.Range("I126").Value = 33337308
$ws.Range("J126").Value = 10599
$ws.Range("K126").Value = 100011924
$ws.Range("L126").Value = 31797
$ws.Range("M126").Value = -100009454
$ws.Range("N126").Value = -36737

$ws.Range("H132").Value = 58826750
$ws.Range("I132").Value = 62503300
$ws.Range("J132").Value = 1995
$ws.Range("K132").Value = 187509900
$ws.Range("L132").Value = 5985
$ws.Range("M132").Value = -187507370
$ws.Range("N132").Value = -11045

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 226169280
$ws.Range("I17").Value = 35000904
$ws.Range("K17").Value = 35000904
$ws.Range("M17").Value = -35000734

$ws.Range("H46").Value = 20000840
$ws.Range("I46").Value = 800.06665
$ws.Range("J46").Value = 50000900
$ws.Range("K46").Value = 800.06665
$ws.Range("L46").Value = 50000900
$ws.Range("M46").Value = -612.06665
$ws.Range("N46").Value = -50001276

$ws.Range("H55").Value = 2428.8147
$ws.Range("I55").Value = 1313.3334
$ws.Range("J55").Value = 3321.2
$ws.Range("K55").Value = 1313.3334
$ws.Range("L55").Value = 3321.2
$ws.Range("M55").Value = -1140.3334
$ws.Range("N55").Value = -3667.2

$ws.Range("H100").Value = 1713.4117
$ws.Range("I100").Value = 905.8
$ws.Range("K100").Value = 905.8
$ws.Range("M100").Value = -364.8

$ws.Range("H104").Value = 46739.2
$ws.Range("J104").Value = 56840
$ws.Range("L104").Value = 56840
$ws.Range("N104").Value = -63828

$ws.Range("H122").Value = 5761.636
$ws.Range("J122").Value = 8064.5
$ws.Range("L122").Value = 24193.5
$ws.Range("N122").Value = -29093.5

$ws.Range("H132").Value = 2962.5557
$ws.Range("I132").Value = 2911.6365
$ws.Range("J132").Value = 3042.5715
$ws.Range("K132").Value = 8734.9095
$ws.Range("L132").Value = 9127.7145
$ws.Range("M132").Value = -6204.9095
$ws.Range("N132").Value = -14187.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 60322.6
$ws.Range("J105").Value = 60322.6
$ws.Range("L105").Value = 60322.6
$ws.Range("N105").Value = -67310.60000000001

$ws.Range("H132").Value = 5114.074
$ws.Range("J132").Value = 11491.444
$ws.Range("L132").Value = 34474.33199999999
$ws.Range("N132").Value = -39534.33199999999

$ws.Range("H136").Value = 18559186
$ws.Range("I136").Value = 26370272
$ws.Range("K136").Value = 79110816
$ws.Range("M136").Value = -79108266
